# Saldo.xlsx edit: apply row additions/removals/replacements described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work bottom-to-top so earlier row numbers stay valid as rows are removed.

# 1) Remove the "Carlos" row (004211922 / Carlos / 3967.63) - originally row 13.
$ws.Rows.Item(13).Delete()

# 2) Remove the four rows Casmarinho/Lais/Patricia/Daura (originally rows 7-10).
$ws.Range("A7:A10").EntireRow.Delete()

# 3) Replace the Luana/Styphany rows (originally rows 4-5) with Fernanda/Walquiria.
$ws.Cells.Item(4,1).NumberFormat = "@"
$ws.Cells.Item(4,1).Value = "000806386"
$ws.Cells.Item(4,2).Value = "Fernanda"
$ws.Cells.Item(4,3).Value = 50515.78

$ws.Cells.Item(5,1).NumberFormat = "@"
$ws.Cells.Item(5,1).Value = "005103059"
$ws.Cells.Item(5,2).Value = "Walquiria"
$ws.Cells.Item(5,3).Value = 32019.92

# 4) Insert a new row 2 for Marcio (008364902 / Marcio / 58839.73), pushing the
#    rest of the table down by one row.
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2,1).NumberFormat = "@"
$ws.Cells.Item(2,1).Value = "008364902"
$ws.Cells.Item(2,2).Value = "Marcio"
$ws.Cells.Item(2,3).Value = 58839.73
